$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing data rows down by one
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Set values for the newly inserted row 2
$ws.Cells.Item(2, 1).Value = 0.6115283966064453
$ws.Cells.Item(2, 2).Value = -0.0169776529073715
$ws.Cells.Item(2, 3).Value = -0.1087901294231414

# Append new rows of data at the end (rows 23-31)
$newRows = @(
    @(-11.22835350036621, -15.27582550048828, -1.236372590065002),
    @(4.94928503036499, -15.6870174407959, 4.060655117034912),
    @(3.654456377029419, -6.942261695861816, 2.724813222885132),
    @(6.618554592132568, 4.984438896179199, -4.156262397766113),
    @(-4.189085960388184, 1.281579732894898, 2.08831787109375),
    @(-3.292665958404541, 1.869073033332825, 3.439073085784912),
    @(-4.698282241821289, 7.67050313949585, -1.46966552734375),
    @(1.039232015609741, 13.21467208862305, -9.619471549987791),
    @(5.204416275024414, -5.015731334686279, -0.8006793856620789)
)

$startRow = 23
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($row, 3).Value = $newRows[$i][2]
}
